# Update the "selected features" worksheet: refresh the feature table with
# the new selection (15 rows of data instead of 8), including new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# First, replicate the formatting of column A (bold + border, style used by
# A2:A8) down through the new rows (A9:A15) so every data row in column A
# keeps the same look once populated.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A9:A15").PasteSpecial(-4122) | Out-Null

# Feature rows: index, feature name, freq, p
$data = @(
    @(2,  8,   "{ankle, hip} (distance) - |temporal| std_min",                               1,    0),
    @(3,  11,  "{ankle, hip} (distance) {diff} - |spectral| rel_pwr_2_to_4_min",              0.9,  0),
    @(4,  33,  "{elbow, shoulder, hip} (angle) {diff} - |spectral| rel_pwr_0.5_to_1_max",     0.9,  0),
    @(5,  38,  "{elbow, shoulder, hip} (angle) {diff} - |spectral| rel_pwr_4_to_6_max",       0.9,  0),
    @(6,  39,  "{elbow, shoulder, hip} (angle) {diff} - |spectral| rel_pwr_4_to_6_min",       0.95, 0),
    @(7,  47,  "{left_ankle, right_ankle} (x_displacement) - |spectral| rel_pwr_0.5_to_1",    0.95, 0),
    @(8,  48,  "{left_ankle, right_ankle} (x_displacement) - |spectral| rel_pwr_1_to_2",      0.8,  0),
    @(9,  55,  "{left_ankle, right_ankle} (x_displacement) {diff} - |spectral| rel_pwr_4_to_6", 0.9, 0),
    @(10, 125, "{nose, middle_shoulder, left_shoulder} (angle) - |spectral| rel_pwr_2_to_4",  0.95, 0),
    @(11, 135, "{pinky, wrist} (distance) {diff} - |spectral| half_pwr_freq",                 0.85, 0),
    @(12, 138, "{pinky, wrist} (distance) {diff} - |spectral| rel_pwr_6_to_128",              1,    0),
    @(13, 151, "{thumb, index} (distance) - |temporal| mad",                                  1,    0),
    @(14, 156, "{thumb, index} (distance) {diff} - |spectral| rel_pwr_2_to_4",                1,    0),
    @(15, 157, "{thumb, index} (distance) {diff} - |spectral| rel_pwr_4_to_6",                0.9,  0)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}
